$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# Insert a new column before column I ("date"), pushing date/legislator_name/legislator_id
# one column to the right, to make room for the new "category" column.
$ws.Columns("I:I").Insert()

# New header cells, in the order they first appear left-to-right / top-to-bottom so
# that new shared strings get appended in the same order as the reference workbook.
$ws.Range("I1").Value = "category"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# Copy header formatting (bold/centered/bordered) from the existing L1 header cell
# onto the two brand-new trailing header cells.
$ws.Range("L1").Copy()
$ws.Range("M1:N1").PasteSpecial(-4122)

# Row 2
$ws.Range("I2").Value = "normal"
$ws.Range("M2").Value = "tmp6101"
$ws.Range("N2").Value = 96

# Row 3
$ws.Range("I3").Value = "normal"
$ws.Range("M3").Value = "tmp6101"
$ws.Range("N3").Value = 97

# Row 4
$ws.Range("I4").Value = "normal"
$ws.Range("M4").Value = "tmp6101"
$ws.Range("N4").Value = 98

$excel.CutCopyMode = 0
